$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at S (shifts nothing right since S was blank/out of range,
# but this makes Excel extend the used range/dimension and inherit formatting
# from the adjacent column to the left, matching the existing column pattern).
$ws.Columns("S:S").Insert()

# Fill in the new 2022 column of data, mirroring the existing year columns.
$ws.Range("S4").Value = 2022
$ws.Range("S5").Value = 0
$ws.Range("S6").Value = 0
$ws.Range("S7").Value = 0
$ws.Range("S8").Value = 0
$ws.Range("S9").Value = 0
$ws.Range("S10").Value = 0
$ws.Range("S11").Value = 0
$ws.Range("S12").Value = 0
$ws.Range("S13").Value = 0
$ws.Range("S14").Value = 0

# Update the sheet's remembered selection, as saved in the workbook.
$ws.Range("R17").Select() | Out-Null
